# Generate Report for Handoff
# Refresh the localization-status report: a new handoff (md) file and new
# target (.xlf) files were generated, along with new handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldId  = "cfaf5365-160b-438b-83e5-dedf244e3a9a"
$newId  = "8f8b09fb-3733-4220-b45c-342dbd04b657"
$oldTok = "3b0a2cf539391c6080af0f74dc4b3ad45c6f5143"
$newTok = "a707f0fd1a74b3fd9620bee3774192b485c7c02b"

$oldMd = "$oldId.md"
$newMd = "$newId.md"

$oldZhXlf = "$oldId.$oldTok.zh-cn.xlf"
$newZhXlf = "$newId.$newTok.zh-cn.xlf"

$oldDeXlf = "$oldId.$oldTok.de-de.xlf"
$newDeXlf = "$newId.$newTok.de-de.xlf"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($hl in $wsOverview.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    }
}
$wsOverview.Range("A2").Value = $newMd
$wsOverview.Range("D2").Value = "2016-03-24 05:05:48"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($hl in $wsZh.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.TextToDisplay -eq $oldZhXlf) {
        $hl.TextToDisplay = $newZhXlf
    }
}
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = "2016-03-24 05:05:44"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($hl in $wsDe.Hyperlinks) {
    if ($hl.TextToDisplay -eq $oldMd) {
        $hl.TextToDisplay = $newMd
    } elseif ($hl.TextToDisplay -eq $oldDeXlf) {
        $hl.TextToDisplay = $newDeXlf
    }
}
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = "2016-03-24 05:05:48"
